$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78 (shifts existing rows 78-119 down to 79-120)
$ws.Range("A78").EntireRow.Insert()

# Populate the newly inserted row 78 with the new weekly record
$ws.Range("A78").Value = 7
$ws.Range("B78").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C78").Value = "Ñuble"
$ws.Range("D78").Value = 44452
$ws.Range("E78").Value = 16
$ws.Range("F78").Value = 100112006
$ws.Range("G78").Value = "Repollo"
$ws.Range("H78").Value = "Crespo record"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 300
$ws.Range("K78").Value = 600
$ws.Range("L78").Value = 650
$ws.Range("M78").Value = 625
$ws.Range("N78").Value = "$/unidad"
$ws.Range("O78").Value = "Provincia de Diguillín"
$ws.Range("P78").Value = 625
$ws.Range("Q78").Value = 1
$ws.Range("R78").Value = "Hortaliza"
